# FluWatch 2022-2023 data refresh: corrected several existing weekly
# counts (wk 38, 41-45) against the published source and appended the
# newly released week-46 row, then left the sheet scrolled/selected
# where the analyst was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Corrections to previously-entered weeks ---
$ws.Range("C5").Value = 35          # week 38

$ws.Range("B8").Value = 120         # week 41
$ws.Range("C8").Value = 110

$ws.Range("B9").Value = 214         # week 42
$ws.Range("C9").Value = 204
$ws.Range("F9").Value = 2.4

$ws.Range("B10").Value = 407        # week 43
$ws.Range("C10").Value = 595
$ws.Range("D10").Value = 23
$ws.Range("F10").Value = 5.6

$ws.Range("B11").Value = 995        # week 44
$ws.Range("C11").Value = 1351
$ws.Range("D11").Value = 54
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 10.9

$ws.Range("B12").Value = 1799       # week 45
$ws.Range("C12").Value = 2265
$ws.Range("D12").Value = 69
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 16

# --- Append the newly published week 46 row ---
$ws.Range("A13").Value = 46
$ws.Range("B13").Value = 3159
$ws.Range("C13").Value = 2213
$ws.Range("D13").Value = 73
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = 19.2
$ws.Range("G13").Value = 0.1

# Carry over the same look (fill/border/font + row height) as the rest
# of the data rows so the new row matches the table formatting.
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows(13).RowHeight = $ws.Rows(12).RowHeight
$excel.CutCopyMode = $false

# --- Leave the view where the analyst ended up: scrolled to the new
# rows, with cell A15 selected ---
$ws.Range("A15").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
